$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 0.1037184923425638
$ws.Range("E4").Value = 0.1523866428317699

$ws.Range("C5").Value = -0.3429736775751868
$ws.Range("E5").Value = -0.1149259972201833

$ws.Range("C6").Value = -0.2064822032187164
$ws.Range("E6").Value = -0.1075365560329056

$ws.Range("C7").Value = 0.005797897989445744
$ws.Range("E7").Value = -0.05125851421730054

$ws.Range("C8").Value = 0.2273278155300318
$ws.Range("E8").Value = 0.001779459905826286

$ws.Range("C9").Value = -0.2723669344146917
$ws.Range("E9").Value = -0.1573879283727764

$ws.Range("C10").Value = -0.05011389829933099
$ws.Range("E10").Value = -0.1100689213476058

$ws.Range("C11").Value = 0.2788833036191596
$ws.Range("E11").Value = -0.1416348838281123

$ws.Range("C12").Value = -0.1459594536071473
$ws.Range("E12").Value = -0.1825768856163368

$ws.Range("C13").Value = -0.2305835819295887
$ws.Range("E13").Value = -0.2087674606261247

$ws.Range("C14").Value = -0.3292346503903532
$ws.Range("E14").Value = -0.09082057608673644

$ws.Range("C15").Value = -0.05652554630790618
$ws.Range("E15").Value = 0.149910086121019

$ws.Range("C16").Value = 1.183936177705625
$ws.Range("E16").Value = 0.1147913851119675

$ws.Range("C17").Value = 0.8096199637471102
$ws.Range("E17").Value = -0.08981539700775309

$ws.Range("C18").Value = -1.085560181261136
$ws.Range("E18").Value = -0.09271976299123352

$ws.Range("C19").Value = 0.5641976336596244
$ws.Range("E19").Value = -0.1178589452312528
